# no-op
